# Auto update photos 2025-11-26 09:46:42
#
# The "Summary" sheet lists product articles with their photo links and
# photo counts. This run's refresh drops every article except "005" -
# i.e. rows 2 ("0026-2"), 3 ("0026-3"), 5 ("009-2") and 6 ("009-4") are
# removed, leaving only the header row plus the "005" row (which shifts
# up to row 2).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete from the bottom up so earlier row numbers stay valid as we go.
$ws.Range("A6:C6").EntireRow.Delete()
$ws.Range("A5:C5").EntireRow.Delete()
$ws.Range("A3:C3").EntireRow.Delete()
$ws.Range("A2:C2").EntireRow.Delete()
